$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 26, shifting existing rows 26-83 down to 27-84.
$ws.Rows("26:26").Insert()

# Populate the newly inserted row 26 with the new record.
# (Most columns mirror the record that used to be in row 26 / now row 27,
#  only D, J, K, L, M, O and P differ per the diff.)
$ws.Range("A26").Value = 10
$ws.Range("B26").Value = "Vega Modelo de Temuco"
$ws.Range("C26").Value = "La Araucanía"
$ws.Range("D26").Value = 44914
$ws.Range("E26").Value = 9
$ws.Range("F26").Value = 100112030
$ws.Range("G26").Value = "Poroto granado"
$ws.Range("H26").Value = "Sin especificar"
$ws.Range("I26").Value = "Primera"
$ws.Range("J26").Value = 110
$ws.Range("K26").Value = 50000
$ws.Range("L26").Value = 50000
$ws.Range("M26").Value = 50000
$ws.Range("N26").Value = "$/saco 25 kilos"
$ws.Range("O26").Value = "Región Metropolitana"
$ws.Range("P26").Value = 2000
$ws.Range("Q26").Value = 25
$ws.Range("R26").Value = "Hortaliza"
